# Apply resume wording/ordering edits described by the commit:
# "Changed some wording for experience and changed order of frameworks"

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Certifications bullet: "Cisco CCNA Netw" + bookmark + "orking for
#    Home and Small " -> single run "Cisco CCNA Networking for Home and
#    Small " (the _GoBack bookmark that used to live here moves down to
#    the "medical data" bullet below).
# ---------------------------------------------------------------------
$r = $d.Content
$search = "Cisco CCNA Netw" + "orking for Home and Small "
$r.Find.Execute($search, $false, $false, $false, $false, $false, $true, 1, `
    $false, "Cisco CCNA Networking for Home and Small ", 2)

# ---------------------------------------------------------------------
# 2) Frameworks bullet: reorder "Django, Ruby on Rails, Backbone" to
#    "Ruby on Rails, Django, Backbone".
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Django, Ruby on Rails, Backbone", $false, $false, $false, `
    $false, $false, $true, 1, $false, "Ruby on Rails, Django, Backbone", 2)

# ---------------------------------------------------------------------
# 3) Experience heading: "UPMC" followed by a run of tab characters ->
#    "University of Pittsburgh Medical Center (UPMC)" followed by fewer
#    tab characters (5 tabs removed from the run of 11).
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("UPMC", $false, $false, $false, $false, $false, $true, 1, `
    $false, "", 0)
$upmcStart = $r.Start
$upmcEnd = $r.End
$delRange = $d.Range($upmcStart, $upmcEnd + 5)
$delRange.Text = "University of Pittsburgh Medical Center (UPMC)"

# ---------------------------------------------------------------------
# 5) "Developed full stack..." bullet rewritten.
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute( `
    "Developed full stack web application to help UPMC better govern and manage all their information.", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Developed full stack web application using Ruby on Rails to help govern and manage medical data.", 2)

# Move the _GoBack bookmark to sit right before the final period of that
# sentence (matches where the diff re-introduces it).
$r = $d.Content
$r.Find.Execute("medical data", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)
$bmPos = $r.End
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------
# 6) "Built data analytics..." bullet rewritten.
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute( `
    "Built data analytics and visualization dashboards using D3 and Highcharts.", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Built data analytics and visualization dashboards using D3 to find pain points and do impact analysis.", 2)

# ---------------------------------------------------------------------
# 7) Fix stray backtick: "...and Gu`lp" -> "...and Gulp".
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Gu" + [char]96 + "lp", $false, $false, $false, $false, `
    $false, $true, 1, $false, "Gulp", 2)
